$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns keep their existing text formatting (values such as
# "591.99" or "8.13" must stay text, matching the workbook's original
# inline-string cell type, rather than being auto-converted to numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.065.47'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.180.48'
$ws.Range('E3').Value = '  -4.16%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.99'
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.37'
$ws.Range('E6').Value = '  -3.94%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.177.77'
$ws.Range('E8').Value = '  -4.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('E10').Value = '  -4.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.28'
$ws.Range('E11').Value = '  -3.90%  '
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.88'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.703.86'
$ws.Range('E15').Value = '  -4.17%  '
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.179.06'
$ws.Range('E17').Value = '  -4.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.049.48'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('E19').Value = '  -3.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.98'
$ws.Range('E20').Value = '  -3.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.95'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.712'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.67'
$ws.Range('E23').Value = '  -6.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.43'
$ws.Range('E24').Value = '  -1.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.50'
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.68'
$ws.Range('E27').Value = '  -2.98%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.81'
$ws.Range('E29').Value = '  -4.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.81'
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('E31').Value = '  -5.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.39'
$ws.Range('E32').Value = '  -5.58%  '
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.38'
$ws.Range('E34').Value = '  -5.91%  '
$ws.Range('E35').Value = '  -6.08%  '
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('E37').Value = '  -1.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0712'
$ws.Range('E38').Value = '  -4.27%  '
$ws.Range('E39').Value = '  -2.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '406.52'
$ws.Range('E40').Value = '  -6.27%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.13'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.67'
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.113'
$ws.Range('E43').Value = '  -7.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.804.22'
$ws.Range('E44').Value = '  -9.09%  '
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.14'
$ws.Range('E46').Value = '  -3.12%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.66'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.06'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.61'
$ws.Range('E50').Value = '  -6.39%  '
$ws.Range('E51').Value = '  -1.86%  '
